$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting columns C:N to D:O.
# Excel's default behaviour copies the formatting of the column to the
# left (B) into the newly inserted column, which is why the new C3 cell
# ends up sharing B3's number-format/border style.
$ws.Columns("C:C").Insert()

# Populate the newly inserted column with the new "1 Durchschnitt" values.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 523.58600000000001

# The column insert also stamps a blank, styled placeholder cell into row 1
# (which carries a row-wide bold format) even though that row never gets a
# value in column C. Rebuild row 1 from scratch so it only contains the
# three cells that actually hold data, matching the original sparse layout.
$ws.Rows(1).Delete()
$ws.Rows(1).Insert()
$ws.Range("A1").Value = "Tabelle"
$ws.Range("B1").Value = "fd_reduced_15.csv"
$ws.Range("F1").Value = "15 columns, 25000 rows"
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

# Match the author's final cell selection.
$ws.Range("G7").Select()
